$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 77
$ws1.Range("F6").Value = 5358
$ws1.Range("F7").Value = 71
$ws1.Range("F8").Value = 883
$ws1.Range("F10").Value = 2375
$ws1.Range("F11").Value = 77
$ws1.Range("F13").Value = 2228

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 94

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 77
$ws4.Range("F6").Value = 5358
$ws4.Range("F7").Value = 94
$ws4.Range("F8").Value = 71
$ws4.Range("F10").Value = 883
$ws4.Range("F12").Value = 2375
$ws4.Range("F13").Value = 77
$ws4.Range("F16").Value = 2228
